$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data (and a few rows whose coin
# order shifted, swapping Coin/Link/Price/Volume while keeping the
# rank index in column A untouched).

$ws.Range("D2").Value = "'28.263.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.34%  "

$ws.Range("D3").Value = "'1.869.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.74%  "

$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'318.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.79%  "

$ws.Range("D6").Value = "'1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("D7").Value = "'0.4383"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.52%  "

$ws.Range("D8").Value = "'0.3697"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.28%  "

$ws.Range("D9").Value = "'0.07510"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.52%  "

$ws.Range("D10").Value = "'0.9379"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.97%  "

$ws.Range("E11").Value = "  -2.71%  "

$ws.Range("D12").Value = "'1.890.10"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.29%  "

$ws.Range("D13").Value = "'6.724"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.91%  "

$ws.Range("D14").Value = "'5.444"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.64%  "

$ws.Range("D15").Value = "'0.06864"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.48%  "

$ws.Range("D16").Value = "'1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.09%  "

$ws.Range("D17").Value = "'82.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.88%  "

$ws.Range("D18").Value = "'0.000009059"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.05%  "

$ws.Range("D19").Value = "'1.002"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.03%  "

$ws.Range("D20").Value = "'15.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.12%  "

$ws.Range("D21").Value = "'28.230.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.49%  "

$ws.Range("D22").Value = "'5.128"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.18%  "

$ws.Range("D23").Value = "'10.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.69%  "

$ws.Range("D24").Value = "'2.125.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.96%  "

$ws.Range("E25").Value = "  -3.17%  "

$ws.Range("D26").Value = "'154.78"
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").Value = "'18.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.02%  "

$ws.Range("D28").Value = "'5.312"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.86%  "

$ws.Range("D29").Value = "'113.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.03%  "

$ws.Range("D30").Value = "'1.727"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.86%  "

$ws.Range("D31").Value = "'0.09057"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.24%  "

$ws.Range("D32").Value = "'0.7991"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.40%  "

$ws.Range("D33").Value = "'4.841"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.69%  "

$ws.Range("D34").Value = "'1.172"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.58%  "

$ws.Range("D35").Value = "'2.925"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.21%  "

$ws.Range("D36").Value = "'1.002"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.06%  "

$ws.Range("D37").Value = "'1.123"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.18%  "

$ws.Range("D38").Value = "'0.05439"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.79%  "

$ws.Range("D39").Value = "'0.01958"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.96%  "

$ws.Range("D40").Value = "'2.908"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.80%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.5253"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.23%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'7.080"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.13%  "

$ws.Range("D43").Value = "'0.1680"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.19%  "

$ws.Range("D44").Value = "'8.744"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.59%  "

$ws.Range("D45").Value = "'0.06746"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.98%  "

$ws.Range("D46").Value = "'0.4879"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.49%  "

$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "'107.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.14%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'1.978"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.69%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'10.49"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.85%  "

$ws.Range("D50").Value = "'0.000002446"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.49%  "

$ws.Range("D51").Value = "'1.678"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.08%  "
